$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(445).Insert()

$ws.Cells.Item(445, 1).Value = 10
$ws.Cells.Item(445, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(445, 3).Value = "La Araucanía"
$ws.Cells.Item(445, 4).Value = 44783
$ws.Cells.Item(445, 5).Value = 9
$ws.Cells.Item(445, 6).Value = 100112023
$ws.Cells.Item(445, 7).Value = "Brócoli"
$ws.Cells.Item(445, 8).Value = "Sin especificar"
$ws.Cells.Item(445, 9).Value = "Primera"
$ws.Cells.Item(445, 10).Value = 1200
$ws.Cells.Item(445, 11).Value = 1000
$ws.Cells.Item(445, 12).Value = 1000
$ws.Cells.Item(445, 13).Value = 1000
$ws.Cells.Item(445, 14).Value = "`$/unidad"
$ws.Cells.Item(445, 15).Value = "Región Metropolitana"
$ws.Cells.Item(445, 16).Value = 1000
$ws.Cells.Item(445, 17).Value = 1
$ws.Cells.Item(445, 18).Value = "Hortaliza"
